$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are written through a temporary "@" (text)
# number format so that numeric-looking strings, e.g. "1.001", are stored
# as text rather than being auto-parsed into a Number -- matching the
# source workbook, where every Price cell is an inline string. The format
# is applied and cleared per cell (instead of via a multi-area union
# range) to sidestep a NumberFormat propagation quirk on union ranges,
# and ClearFormats() afterwards removes the temporary format so the cell
# keeps its original (default) style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.939.14'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.92'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.91'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4693'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3701'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07391'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8729'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.51'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.35'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.372'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.01'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07078'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.522'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008731'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.960.07'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.334'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.032.04'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.904'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.45'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.210'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.42'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.334'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.86'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08944'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7713'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.169'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.507'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.901'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("E37").Value = '  -2.79%  '
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05298'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.344'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.945'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5360'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.377'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1674'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.479'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4974'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.45'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.678'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.000'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.30'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("E51").Value = '  +0.30%  '

Write-Host "Updated cryptos list"